$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($CellRef, $Text) {
    $cell = $ws.Range($CellRef)
    # Leading apostrophe forces Excel to store the numeric-looking
    # string as literal text instead of coercing it to a number.
    $cell.Value = "'" + $Text
    # Re-apply the default "Normal" style so the quote-prefix flag
    # added by the apostrophe trick does not linger on the cell -
    # keeps the cell style identical to before the edit.
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '311.46'
Set-TextCell 'E2' '-0.12%'
Set-TextCell 'D3' '37.69'
Set-TextCell 'E3' '-1.71%'
Set-TextCell 'D4' '5.080'
Set-TextCell 'E4' '-0.90%'
Set-TextCell 'D5' '0.07779'
Set-TextCell 'E5' '-4.03%'
Set-TextCell 'D6' '4.354'
Set-TextCell 'E6' '-2.76%'
Set-TextCell 'D7' '8.226'
Set-TextCell 'E7' '-1.08%'
Set-TextCell 'D8' '1.884'
Set-TextCell 'E8' '-3.91%'
Set-TextCell 'D9' '2.850'
Set-TextCell 'E9' '-10.72%'
Set-TextCell 'D10' '0.9201'
Set-TextCell 'E10' '-2.18%'
Set-TextCell 'D11' '0.1195'
Set-TextCell 'E11' '-9.92%'
Set-TextCell 'E12' '-1.67%'
Set-TextCell 'D13' '0.09336'
Set-TextCell 'E13' '3.67%'
Set-TextCell 'D14' '0.03397'
Set-TextCell 'E14' '-2.62%'
Set-TextCell 'D15' '0.09679'
Set-TextCell 'E15' '-0.30%'
Set-TextCell 'D16' '0.001383'
Set-TextCell 'E16' '-1.81%'
Set-TextCell 'D17' '0.005753'
Set-TextCell 'E17' '-4.89%'
Set-TextCell 'E18' '-0.40%'
Set-TextCell 'D19' '0.3404'
Set-TextCell 'E19' '-1.80%'
Set-TextCell 'D20' '5.267'
Set-TextCell 'E20' '4.83%'
Set-TextCell 'D21' '0.1282'
Set-TextCell 'E21' '-0.65%'
Set-TextCell 'D22' '0.2589'
Set-TextCell 'E22' '3.89%'
Set-TextCell 'E23' '5,585.19%'
Set-TextCell 'D24' '0.04373'
Set-TextCell 'E24' '-0.07%'
Set-TextCell 'D25' '0.001214'
Set-TextCell 'E25' '-2.56%'
Set-TextCell 'D26' '0.004260'
Set-TextCell 'E26' '-9.99%'
Set-TextCell 'E27' '-66.74%'
Set-TextCell 'E39' '-3.59%'
Set-TextCell 'D40' '0.04981'
Set-TextCell 'E40' '-4.93%'
Set-TextCell 'D41' '0.007672'
Set-TextCell 'E41' '0.54%'
Set-TextCell 'D42' '0.009901'
Set-TextCell 'E42' '-3.02%'
Set-TextCell 'D43' '0.1346'
Set-TextCell 'E43' '-3.02%'
Set-TextCell 'D44' '0.002060'
Set-TextCell 'E44' '1.09%'
Set-TextCell 'D45' '0.008813'
Set-TextCell 'E45' '-3.39%'
Set-TextCell 'D46' '0.00006660'
Set-TextCell 'E46' '0.61%'
Set-TextCell 'E47' '-0.45%'
Set-TextCell 'D48' '0.002912'
Set-TextCell 'E48' '-3.54%'
Set-TextCell 'E50' '-0.45%'
Set-TextCell 'E51' '-0.45%'
